# Applies the "Add new customer, session, booking all working.
# Pre-double book prevention trial" edit to the Project Planning document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Still To Do:" checklist (numId=6 list)
#    Old:
#      Wireframe diagrams
#      All CSS
#      Customer DB / New / Update / Delete
#      Session DB / New
#      Booking DB / New / (delete) / (update)
#      Update and complete diagrams
#      Evidence on testing
#    New:
#      All CSS
#      Booking DB
#       (delete)
#      (update)
#      Update and complete diagrams
#      Evidence on testing
# ---------------------------------------------------------------

# 1a) "Wireframe diagrams" -> "All CSS"
$pWire = $d.Paragraphs.Item(49)
$rWire = $d.Range($pWire.Range.Start, $pWire.Range.End - 1)
$rWire.Text = "All CSS"

# 1b) Remove the now-superseded "All CSS" / Customer DB (+ New/Update/Delete) /
#     Session DB (+ New) paragraphs - 7 paragraphs in total, directly after the
#     paragraph just edited, leaving "Booking DB" immediately following.
$pDelStart = $d.Paragraphs.Item(50)
$pDelEnd = $d.Paragraphs.Item(56)
$delRange = $d.Range($pDelStart.Range.Start, $pDelEnd.Range.End)
$delRange.Delete()

# 1c) The sub-bullet "New" (under "Booking DB") becomes a leading space, and the
#     following "(delete)" paragraph is folded into it as a second run (its own
#     paragraph mark is removed), yielding one bullet reading " (delete)".
$pBookingNew = $d.Paragraphs.Item(51)
$rBookingNew = $d.Range($pBookingNew.Range.Start, $pBookingNew.Range.End - 1)
$rBookingNew.Text = " "
$pBookingNew2 = $d.Paragraphs.Item(51)
$markRange = $d.Range($pBookingNew2.Range.End - 1, $pBookingNew2.Range.End)
$markRange.Delete()

# ---------------------------------------------------------------
# 2) "Extensions:" list (numId=7) - insert a new bullet right before
#    "Customer Class".
# ---------------------------------------------------------------
$pCustClass = $d.Paragraphs.Item(57)
$pCustClass.Range.InsertParagraphBefore()
$pNewBullet = $d.Paragraphs.Item(57)
$rNewBullet = $d.Range($pNewBullet.Range.Start, $pNewBullet.Range.End - 1)
$rNewBullet.Text = "If a customer has already been booked into the class, they can no longer attend that same class"

# ---------------------------------------------------------------
# 3) Turn the blank paragraph right after "...a capacity" into a new
#    Extensions bullet ("Change all date/times to actual date/times"),
#    and move the hidden "_GoBack" bookmark here (it previously sat on
#    the "Edit function likewise..." paragraph further down).
# ---------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(68)
$pBlank.Style = "List Paragraph"
$pBlank.Range.ListFormat.ListLevelNumber = 1

$rBlank = $d.Paragraphs.Item(68).Range
$rBlankNoMark = $d.Range($rBlank.Start, $rBlank.End - 1)
$rBlankNoMark.Text = "Change all date/times to actual date/times"

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# A bookmark collapsed exactly at a paragraph boundary gets mis-anchored by
# this host, so temporarily append a placeholder run, anchor the (now
# mid-paragraph) collapsed bookmark against that boundary, then remove the
# placeholder - leaving the bookmark sitting right after the real text.
$pFinal = $d.Paragraphs.Item(68)
$rFinalNoMark = $d.Range($pFinal.Range.Start, $pFinal.Range.End - 1)
$rFinalNoMark.InsertAfter("Z")

$pFinal2 = $d.Paragraphs.Item(68)
$boundaryPos = $pFinal2.Range.End - 2
$midRange = $d.Range($boundaryPos, $boundaryPos)
$d.Bookmarks.Add("_GoBack", $midRange)

$pFinal3 = $d.Paragraphs.Item(68)
$placeholderRange = $d.Range($pFinal3.Range.End - 2, $pFinal3.Range.End - 1)
$placeholderRange.Delete()
